$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.1800000000005
$ws.Range("H2").Value = 0.0000000000000003929993007522678
$ws.Range("K2").Value = 57.06381856368852
$ws.Range("L2").Value = "[45.277677653523924, 68.84995947385312]"
$ws.Range("O2").Value = 1.415131825941348
$ws.Range("P2").Value = "[1.1887107337907317, 1.6415529180919641]"
$ws.Range("S2").Value = 63.97707319730128
$ws.Range("T2").Value = "[56.67522885135239, 71.27891754325017]"
$ws.Range("W2").Value = 19.50882882882922
$ws.Range("X2").Value = 18.60144144144181
$ws.Range("Y2").Value = 20.41621621621663

# Row 3 updates
$ws.Range("E3").Value = 23.08000000000017
$ws.Range("G3").Value = 0.000000000000001887379141862766
$ws.Range("H3").Value = 0.000000000000005450914489134343
$ws.Range("K3").Value = 53.39158391771458
$ws.Range("L3").Value = "[38.102977960833186, 68.68018987459597]"
$ws.Range("M3").Value = 0.0000000001586959452737347
$ws.Range("N3").Value = 0.0000000001586959452737347
$ws.Range("O3").Value = -0.8427896207828471
$ws.Range("P3").Value = "[-1.1321054607530776, -0.5534737808126167]"
$ws.Range("Q3").Value = 0.00000005091722310446301
$ws.Range("R3").Value = 0.00000005091722310446301
$ws.Range("S3").Value = 60.81332351574046
$ws.Range("T3").Value = "[52.73589399297968, 68.89075303850123]"
$ws.Range("W3").Value = 3.09581581581584
$ws.Range("X3").Value = 2.033073073073092
$ws.Range("Y3").Value = 4.158558558558587
